# Update leve profit calculation cells per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1803.5869
$ws.Range("J17").Value = 1877.1951
$ws.Range("L17").Value = 5631.5853
$ws.Range("N17").Value = -5967.5853
$ws.Range("H76").Value = 8428.429
$ws.Range("I76").Value = 9499.5
$ws.Range("J76").Value = 8000
$ws.Range("K76").Value = 9499.5
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = -9184.5
$ws.Range("N76").Value = -8630
$ws.Range("H79").Value = 8428.429
$ws.Range("I79").Value = 9499.5
$ws.Range("J79").Value = 8000
$ws.Range("K79").Value = 9499.5
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = -8407.5
$ws.Range("N79").Value = -10184
$ws.Range("H87").Value = 64887.25
$ws.Range("J87").Value = 94774.5
$ws.Range("L87").Value = 94774.5
$ws.Range("N87").Value = -97270.5
$ws.Range("H90").Value = 64887.25
$ws.Range("J90").Value = 94774.5
$ws.Range("L90").Value = 284323.5
$ws.Range("N90").Value = -296803.5
$ws.Range("H112").Value = 1687.3112
$ws.Range("J112").Value = 1746.9048
$ws.Range("L112").Value = 5240.7144
$ws.Range("N112").Value = -7456.7144
$ws.Range("H113").Value = 3967.125
$ws.Range("I113").Value = 2589.6667
$ws.Range("J113").Value = 8099.5
$ws.Range("K113").Value = 2589.6667
$ws.Range("L113").Value = 8099.5
$ws.Range("M113").Value = 664.3332999999998
$ws.Range("N113").Value = -14607.5
$ws.Range("H137").Value = 845153.4399999999
$ws.Range("I137").Value = 1215.0588
$ws.Range("K137").Value = 3645.1764
$ws.Range("M137").Value = -1095.1764
$ws.Range("H138").Value = 4398.971
$ws.Range("I138").Value = 2855.4614
$ws.Range("J138").Value = 5311.0454
$ws.Range("K138").Value = 8566.3842
$ws.Range("L138").Value = 15933.1362
$ws.Range("M138").Value = -3426.3842
$ws.Range("N138").Value = -26213.1362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2101.9092
$ws.Range("I2").Value = 1913.7931
$ws.Range("K2").Value = 1913.7931
$ws.Range("M2").Value = -1800.7931
$ws.Range("H74").Value = 3032.366
$ws.Range("I74").Value = 2234.9119
$ws.Range("K74").Value = 2234.9119
$ws.Range("M74").Value = -1360.9119
$ws.Range("H77").Value = 3032.366
$ws.Range("I77").Value = 2234.9119
$ws.Range("K77").Value = 11174.5595
$ws.Range("M77").Value = -6806.559499999999
$ws.Range("H116").Value = 2101.9092
$ws.Range("I116").Value = 1913.7931
$ws.Range("K116").Value = 1913.7931
$ws.Range("M116").Value = 380.2068999999999
$ws.Range("H132").Value = 7249.0835
$ws.Range("I132").Value = 4670.7144
$ws.Range("K132").Value = 14012.1432
$ws.Range("M132").Value = -11482.1432
$ws.Range("H134").Value = 87723.336
$ws.Range("J134").Value = 87723.336
$ws.Range("L134").Value = 87723.336
$ws.Range("N134").Value = -97863.336
$ws.Range("H139").Value = 103237.664
$ws.Range("J139").Value = 103237.664
$ws.Range("L139").Value = 103237.664
$ws.Range("N139").Value = -113517.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2101.9092
$ws.Range("I3").Value = 1913.7931
$ws.Range("K3").Value = 1913.7931
$ws.Range("M3").Value = -1799.7931
$ws.Range("H20").Value = 2450.2666
$ws.Range("I20").Value = 2543.6667
$ws.Range("J20").Value = 2310.1667
$ws.Range("K20").Value = 2543.6667
$ws.Range("L20").Value = 2310.1667
$ws.Range("M20").Value = -2296.6667
$ws.Range("N20").Value = -2804.1667
$ws.Range("H96").Value = 13299.667
$ws.Range("I96").Value = 13299.667
$ws.Range("K96").Value = 13299.667
$ws.Range("M96").Value = -10553.667
$ws.Range("H99").Value = 22055.357
$ws.Range("I99").Value = 30424.334
$ws.Range("K99").Value = 30424.334
$ws.Range("M99").Value = -28926.334
$ws.Range("H135").Value = 99980
$ws.Range("J135").Value = 99980
$ws.Range("L135").Value = 99980
$ws.Range("N135").Value = -110120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1047
$ws.Range("I105").Value = 1047
$ws.Range("K105").Value = 1047
$ws.Range("M105").Value = 700
$ws.Range("H140").Value = 79497
$ws.Range("J140").Value = 79497
$ws.Range("L140").Value = 79497
$ws.Range("N140").Value = -89857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1117.409
$ws.Range("I92").Value = 916.4666999999999
$ws.Range("K92").Value = 2749.4001
$ws.Range("M92").Value = -1501.4001
$ws.Range("H113").Value = 2164.4443
$ws.Range("J113").Value = 2279.1667
$ws.Range("L113").Value = 6837.500100000001
$ws.Range("N113").Value = -11177.5001
$ws.Range("H121").Value = 10527483
$ws.Range("I121").Value = 40000292
$ws.Range("J121").Value = 1479.4286
$ws.Range("K121").Value = 120000876
$ws.Range("L121").Value = 4438.2858
$ws.Range("M121").Value = -119999566
$ws.Range("N121").Value = -7058.2858
$ws.Range("H131").Value = 12321.429
$ws.Range("J131").Value = 14264.444
$ws.Range("L131").Value = 42793.33199999999
$ws.Range("N131").Value = -52873.33199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1121364.4
$ws.Range("I21").Value = 2003295.8
$ws.Range("K21").Value = 2003295.8
$ws.Range("M21").Value = -2003122.8
$ws.Range("H30").Value = 1121364.4
$ws.Range("I30").Value = 2003295.8
$ws.Range("K30").Value = 2003295.8
$ws.Range("M30").Value = -2003190.8
$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30368
$ws.Range("H70").Value = 4333.3335
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -4730
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 4333.3335
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -4064
$ws.Range("N73").Value = -5872
$ws.Range("H80").Value = 111922.086
$ws.Range("I80").Value = 213798.67
$ws.Range("K80").Value = 213798.67
$ws.Range("M80").Value = -212800.67
$ws.Range("H83").Value = 111922.086
$ws.Range("I83").Value = 213798.67
$ws.Range("K83").Value = 1068993.35
$ws.Range("M83").Value = -1064001.35
$ws.Range("H132").Value = 1007332.5
$ws.Range("I132").Value = 1723172.8
$ws.Range("K132").Value = 5169518.4
$ws.Range("M132").Value = -5166988.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2526.3333
$ws.Range("I7").Value = 2517.2727
$ws.Range("K7").Value = 2517.2727
$ws.Range("M7").Value = -2405.2727
$ws.Range("H40").Value = 5916.1113
$ws.Range("I40").Value = 5155.6875
$ws.Range("K40").Value = 5155.6875
$ws.Range("M40").Value = -5019.6875
$ws.Range("H61").Value = 3759.625
$ws.Range("I61").Value = 2110.1428
$ws.Range("K61").Value = 2110.1428
$ws.Range("M61").Value = -1908.1428
$ws.Range("H113").Value = 3759.625
$ws.Range("I113").Value = 2110.1428
$ws.Range("K113").Value = 2110.1428
$ws.Range("M113").Value = 59.85719999999992
$ws.Range("H126").Value = 2526.3333
$ws.Range("I126").Value = 2517.2727
$ws.Range("K126").Value = 7551.8181
$ws.Range("M126").Value = -5081.8181
$ws.Range("H132").Value = 826207.4
$ws.Range("I132").Value = 962950.25
$ws.Range("K132").Value = 2888850.75
$ws.Range("M132").Value = -2886320.75
$ws.Range("H141").Value = 79350
$ws.Range("J141").Value = 79350
$ws.Range("L141").Value = 79350
$ws.Range("N141").Value = -89710

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 58781.75
$ws.Range("J68").Value = 55712
$ws.Range("L68").Value = 55712
$ws.Range("N68").Value = -57334
$ws.Range("H70").Value = 50105
$ws.Range("J70").Value = 50105
$ws.Range("L70").Value = 50105
$ws.Range("N70").Value = -50735
$ws.Range("H71").Value = 58781.75
$ws.Range("J71").Value = 55712
$ws.Range("L71").Value = 167136
$ws.Range("N71").Value = -175248
$ws.Range("H73").Value = 50105
$ws.Range("J73").Value = 50105
$ws.Range("L73").Value = 50105
$ws.Range("N73").Value = -52289
$ws.Range("H111").Value = 65000
$ws.Range("J111").Value = 65000
$ws.Range("L111").Value = 65000
$ws.Range("N111").Value = -73180
$ws.Range("H126").Value = 3330.2856
$ws.Range("I126").Value = 2817.2307
$ws.Range("K126").Value = 8451.6921
$ws.Range("M126").Value = -5981.6921
